$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C: short "ok to miss" annotations against the powers-of-2 rows
# (rows 1,2,3,6,10 correspond to A values 1,2,4,8,16), plus an explanatory
# comment on row 1 and an ellipsis placeholder on row 15 (A=32).
$ws.Range("C2").Value = "ok to miss"
$ws.Range("C15").Value = "…"
$ws.Range("C1").Value = "ok to miss as we already check for powers of 2 and overlay a superset of on bits"
$ws.Range("C3").Value = "ok to miss"
$ws.Range("C6").Value = "ok to miss"
$ws.Range("C10").Value = "ok to miss"

# Give column C a sensible custom width (closest reachable value to 10.71 chars).
$ws.Columns.Item(3).ColumnWidth = 9.83

# Leave the selection where the editing left off.
$ws.Range("E17").Select() | Out-Null
